# appPackage removed from signIn and checkerSignIn sheets
# (and the "devices" sheet's iOS sample row swapped for an Android one)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# devices sheet: the sample row goes from an iOS device to an Android one.
# B2 (platformName) iOS -> Android, M2 (appPackage) gets the checker app id,
# and the now-unused model/automationName/bundleId-value cells are cleared.
# ---------------------------------------------------------------------------
$wsDevices = $wb.Worksheets.Item("devices")
$wsDevices.Range("B2").Value = "Android"
$wsDevices.Range("E2").ClearContents()
$wsDevices.Range("K2").ClearContents()
$wsDevices.Range("L2").ClearContents()
$wsDevices.Range("M2").Value = "au.gov.nsw.onegov.app.checker.uat"

# ---------------------------------------------------------------------------
# signIn sheet: drop the appName column (K) entirely, shifting the
# postal_Address / lic_OwnerName columns one to the left.
# ---------------------------------------------------------------------------
$wsSignIn = $wb.Worksheets.Item("signIn")
$wsSignIn.Columns("K").Delete()

# ---------------------------------------------------------------------------
# checkerSignIn sheet: drop the appName column (D) entirely, shifting the
# licenceNo..Address columns one to the left.
# ---------------------------------------------------------------------------
$wsChecker = $wb.Worksheets.Item("checkerSignIn")
$wsChecker.Columns("D").Delete()

Write-Output "done"
